# Eway bill Unique Constraints Added.
#
# Business changes applied to Sheet1:
#  - EwayBillNo column (B2:B5) switches from a numeric placeholder
#    (987654321123) to a textual "unique constraint" placeholder value
#    "111111111111" (stored as text, not a number).
#  - The ApprovalToSend/ApprovalToReceive flags (P, Q) and the
#    ReceiverRemark (S) column are reset from numeric flags (1 / 0) or a
#    timestamp to the placeholder text "-".
#  - Status (T) moves from "Rejected" to "Pending".
#  - A new (mostly empty) row 7 is added with S7 = "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the cells that need to end up holding genuine text even though some
# of the new values ("111111111111") look numeric. Flipping NumberFormat
# to "@" (Text) first forces Excel to keep the literal characters instead
# of parsing them into a number; doing it once across the whole multi-area
# range (instead of per cell) keeps a single shared style definition
# instead of one per cell. The style is set back to "Normal" afterwards so
# no cell is left with an explicit/visible style difference.
$textRange = $ws.Range("B2:B5,P2:Q5,S2:S5,S7")
$textRange.NumberFormat = "@"

# --- EwayBillNo placeholder: numeric -> text "111111111111" ---
$ws.Range("B2").Value = "111111111111"
$ws.Range("B3").Value = "111111111111"
$ws.Range("B4").Value = "111111111111"
$ws.Range("B5").Value = "111111111111"

# --- ApprovalToSend / ApprovalToReceive / ReceiverRemark -> "-" ---
foreach ($row in 2..5) {
    $ws.Range("P$row").Value = "-"
    $ws.Range("Q$row").Value = "-"
    $ws.Range("S$row").Value = "-"
}

# --- New row 7: just a "-" placeholder in ReceiverRemark ---
$ws.Range("S7").Value = "-"

# Drop the temporary Text formatting now that every value is safely stored.
$textRange.Style = "Normal"

# --- Status: Rejected -> Pending ---
foreach ($row in 2..5) {
    $ws.Range("T$row").Value = "Pending"
}

# --- Selection bookkeeping (cosmetic: matches the saved sheetView) ---
$ws.Range("P2:Q5").Select()
